$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $styleDonor, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $ws.Range($styleDonor).Style
}

Set-TextCell $ws "D2" "B2" "26.334.02"
Set-TextCell $ws "E2" "B2" "  -3.52%  "
Set-TextCell $ws "D3" "B3" "1.665.45"
Set-TextCell $ws "E3" "B3" "  -2.55%  "
Set-TextCell $ws "E4" "B4" "  +0.50%  "
Set-TextCell $ws "D5" "B5" "218.81"
Set-TextCell $ws "E5" "B5" "  -2.32%  "
Set-TextCell $ws "D6" "B6" "0.5159"
Set-TextCell $ws "E6" "B6" "  -3.30%  "
Set-TextCell $ws "D7" "B7" "1.009"
Set-TextCell $ws "E7" "B7" "  +0.44%  "
Set-TextCell $ws "D8" "B8" "0.06439"
Set-TextCell $ws "E8" "B8" "  -2.26%  "
Set-TextCell $ws "D9" "B9" "0.2565"
Set-TextCell $ws "E9" "B9" "  -3.57%  "
Set-TextCell $ws "D10" "B10" "19.95"
Set-TextCell $ws "E10" "B10" "  -4.11%  "
Set-TextCell $ws "D11" "B11" "0.07667"
Set-TextCell $ws "E11" "B11" "  +0.59%  "
Set-TextCell $ws "D12" "B12" "1.684.80"
Set-TextCell $ws "E12" "B12" "  -1.62%  "
Set-TextCell $ws "D13" "B13" "4.334"
Set-TextCell $ws "E13" "B13" "  -5.11%  "
Set-TextCell $ws "D14" "B14" "1.898.22"
Set-TextCell $ws "E14" "B14" "  -2.38%  "
Set-TextCell $ws "D15" "B15" "0.5541"
Set-TextCell $ws "E15" "B15" "  -3.32%  "
Set-TextCell $ws "D16" "B16" "0.0₅8043"
Set-TextCell $ws "E16" "B16" "  -1.50%  "
Set-TextCell $ws "D17" "B17" "64.61"
Set-TextCell $ws "E17" "B17" "  -4.55%  "
Set-TextCell $ws "D18" "B18" "26.379.40"
Set-TextCell $ws "E18" "B18" "  -3.42%  "
Set-TextCell $ws "D19" "B19" "1.006"
Set-TextCell $ws "E19" "B19" "  +0.27%  "
Set-TextCell $ws "D20" "B20" "210.28"
Set-TextCell $ws "E20" "B20" "  -2.61%  "
Set-TextCell $ws "D21" "B21" "4.408"
Set-TextCell $ws "E21" "B21" "  -5.51%  "
Set-TextCell $ws "D22" "B22" "10.10"
Set-TextCell $ws "E22" "B22" "  -3.29%  "
Set-TextCell $ws "D23" "B23" "5.888"
Set-TextCell $ws "E23" "B23" "  -1.33%  "
Set-TextCell $ws "D24" "B24" "1.009"
Set-TextCell $ws "E24" "B24" "  +0.35%  "
Set-TextCell $ws "D25" "B25" "145.24"
Set-TextCell $ws "E25" "B25" "  +2.29%  "
Set-TextCell $ws "D26" "B26" "1.736"
Set-TextCell $ws "E26" "B26" "  -1.27%  "
Set-TextCell $ws "D27" "B27" "0.1165"
Set-TextCell $ws "E27" "B27" "  -4.02%  "
Set-TextCell $ws "D28" "B28" "6.994"
Set-TextCell $ws "E28" "B28" "  -3.82%  "
Set-TextCell $ws "D29" "B29" "15.80"
Set-TextCell $ws "E29" "B29" "  -3.16%  "
Set-TextCell $ws "D30" "B30" "0.05250"
Set-TextCell $ws "E30" "B30" "  -2.81%  "
Set-TextCell $ws "D31" "B31" "1.262"
Set-TextCell $ws "E31" "B31" "  -2.50%  "
Set-TextCell $ws "D32" "B32" "3.369"
Set-TextCell $ws "E32" "B32" "  -3.67%  "
Set-TextCell $ws "D33" "B33" "3.217"
Set-TextCell $ws "E33" "B33" "  -6.03%  "
Set-TextCell $ws "D34" "B34" "1.570"
Set-TextCell $ws "E34" "B34" "  -4.47%  "
Set-TextCell $ws "D35" "B35" "2.757"
Set-TextCell $ws "E35" "B35" "  -4.19%  "
Set-TextCell $ws "D36" "B36" "2.379"
Set-TextCell $ws "E36" "B36" "  -1.33%  "
Set-TextCell $ws "D37" "B37" "0.9270"
Set-TextCell $ws "E37" "B37" "  -2.32%  "
Set-TextCell $ws "D38" "B38" "0.5719"
Set-TextCell $ws "E38" "B38" "  -2.41%  "
Set-TextCell $ws "D39" "B39" "1.148.49"
Set-TextCell $ws "E39" "B39" "  +9.94%  "
Set-TextCell $ws "D40" "B40" "0.01600"
Set-TextCell $ws "E40" "B40" "  -1.76%  "
Set-TextCell $ws "D41" "B41" "0.8477"
Set-TextCell $ws "E41" "B41" "  +0.60%  "
Set-TextCell $ws "D42" "B42" "1.008"
Set-TextCell $ws "E42" "B42" "  +0.29%  "
Set-TextCell $ws "D43" "B43" "5.651"
Set-TextCell $ws "E43" "B43" "  -3.64%  "
Set-TextCell $ws "D44" "B44" "100.06"
Set-TextCell $ws "E44" "B44" "  -0.77%  "
Set-TextCell $ws "D45" "B45" "1.806.76"
Set-TextCell $ws "E45" "B45" "  -2.43%  "
Set-TextCell $ws "E46" "B46" "  +0.00%  "
Set-TextCell $ws "D47" "B47" "0.4499"
Set-TextCell $ws "E47" "B47" "  -0.12%  "
Set-TextCell $ws "D48" "B48" "55.98"
Set-TextCell $ws "E48" "B48" "  -3.50%  "
Set-TextCell $ws "D49" "B49" "1.007"
Set-TextCell $ws "E49" "B49" "  +0.34%  "
Set-TextCell $ws "D50" "B50" "7.910"
Set-TextCell $ws "E50" "B50" "  -1.87%  "
Set-TextCell $ws "D51" "B51" "0.05107"
Set-TextCell $ws "E51" "B51" "  -2.55%  "
